# edit.ps1 - apply the changes described by the diff to the active document.
#
# Summary of the change:
#   1. Remove the two leading empty paragraphs at the very top of the
#      document (they only carried pPr/rPr formatting, no text).
#   2. The paragraph "___________№______" is removed, and the following
#      paragraph ("на №{nStud} от {cDate}") is rewritten in place to
#      "{cDate}№ {nStud} /139-20п" with the exact run-level formatting
#      (including the spell-check proofErr bookmarks around the merge
#      field names) shown in the target markup.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the two leading empty paragraphs.
# ---------------------------------------------------------------------------
$NUM_LEADING_EMPTY = 2
for ($n = 0; $n -lt $NUM_LEADING_EMPTY; $n++) {
    $lead = $d.Paragraphs.Item(1)
    if ($lead.Range.Text.Length -le 1) {
        $lead.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Find the "___________№______" paragraph, delete it (merging it away),
#    then overwrite the paragraph that takes its place with the new
#    "{cDate}№ {nStud} /139-20п" content.
# ---------------------------------------------------------------------------
$markerText = "___________" + [char]8470 + "______"
$numParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith($markerText)) {
        $numParaIndex = $i
        break
    }
}

if ($numParaIndex -gt 0) {
    $numPara = $d.Paragraphs.Item($numParaIndex)
    $numPara.Range.Delete()

    $mergedPara = $d.Paragraphs.Item($numParaIndex)

    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>{</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>cDate</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>}</w:t></w:r>' +
        '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">№ </w:t></w:r>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>{</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>nStud</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>}</w:t></w:r>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>/139-20</w:t></w:r>' +
        '<w:r><w:rPr><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>п</w:t></w:r>' +
        '</w:p>'

    $insertResult = $mergedPara.Range.InsertXML($newParaXml)
}
